$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38 previously described R1,R5 as R0402 / 39 / C25110 / 39Ohm desc.
# Change it to match the R0603 / 1.0k / C21190 / 1kOhm 0603 entry (same
# part as used for R45,R46,R62,R63,R85,R86,R87 on row 42).
$ws.Range("A38").Value = "R0603"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "1.0k"
$ws.Range("D38").Value = "C21190"
$ws.Range("E38").Value = "1/10W Thick Film Resistors 75V ±1% ±100ppm/℃ -55℃~+155℃ 1kΩ 0603  Chip Resistor - Surface Mount ROHS"
